# AICTE_STEGANOGRAPHY.pptx edit:
#   Slide 9 ("GitHub Link"), shape "Content Placeholder 2":
#     - fix/typo the repo URL text to "https://github.com/sSwathi-2003/STEGANOGRAPHY.git"
#     - turn the whole URL into a hyperlink pointing at that same URL
#     - leave a trailing empty (no-bullet) paragraph after it
#
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)

$url = "https://github.com/sSwathi-2003/STEGANOGRAPHY.git"

# Rewrite the paragraph text (typo'd URL) and add a second, empty paragraph
# after it in a single assignment so the paragraph break is preserved.
$tr = $sh.TextFrame.TextRange
$tr.Text = $url + "`r"

# Turn the URL text into a hyperlink. Doing this over four matching
# sub-ranges mirrors how PowerPoint itself breaks the run up around the
# in-place spelling edit, producing the same run boundaries:
#   "https://github." / "com/sSwathi-2003" / "/STEGANOGRAPHY" / ".git"
$tr2 = $sh.TextFrame.TextRange
$r1 = $tr2.Characters(1, 15)
$r2 = $tr2.Characters(16, 16)
$r3 = $tr2.Characters(32, 14)
$r4 = $tr2.Characters(46, 4)
$r1.ActionSettings.Item(1).Hyperlink.Address = $url
$r2.ActionSettings.Item(1).Hyperlink.Address = $url
$r3.ActionSettings.Item(1).Hyperlink.Address = $url
$r4.ActionSettings.Item(1).Hyperlink.Address = $url

# The trailing, empty second paragraph has its bullet switched off.
$tf2 = $sh.TextFrame2
$paras = $tf2.TextRange
$para2 = $paras.Paragraphs(2, 1)
$para2.ParagraphFormat.Bullet.Type = 0
